$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ADC pin assignment updated: re-order which current-sense net/description
# pairs line up with pins 10-13 (rows 12-15) and pins 18-19 (rows 20-21) in
# columns G (net name) and I (description).
$ws.Range("G12").Value = "A3V3_SENSE"
$ws.Range("I12").Value = "A3V3 rail current sense"

$ws.Range("G13").Value = "A1V5_SENSE"
$ws.Range("I13").Value = "A1V5 rail current sense"

$ws.Range("G14").Value = "VOUT_SENSE"
$ws.Range("I14").Value = "VOUT rail current sense"

$ws.Range("G15").Value = "D1V5_SENSE"
$ws.Range("I15").Value = "D1V5 rail current sense"

$ws.Range("G20").Value = "DXVY_SENSE"
$ws.Range("I20").Value = "DXVY rail current sense"

$ws.Range("G21").Value = "D3V3_SENSE"
$ws.Range("I21").Value = "D3V3 rail current sense"

# Leave the cursor where the author left it after making the edit.
$ws.Range("I12").Select()
